$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1806205
$ws.Range("H2").Value = 0.361241
$ws.Range("I2").Value = 0.7284509268949775
$ws.Range("J2").Value = 0.7284509268949775
$ws.Range("M2").Value = 5.380673
$ws.Range("N2").Value = 10.761346
$ws.Range("O2").Value = 0.1901127853194472
$ws.Range("P2").Value = 0.1459548742817064
$ws.Range("Q2").Value = 0.9718598475964999
$ws.Range("R2").Value = 3.887439390386
$ws.Range("S2").Value = 0.1384878346805372
$ws.Range("T2").Value = 0.1063209634553489
$ws.Range("G3").Value = 0.1806205
$ws.Range("H3").Value = 0.361241
$ws.Range("I3").Value = 0.7284509268949775
$ws.Range("J3").Value = 0.7284509268949775
$ws.Range("M3").Value = 5.405099333333332
$ws.Range("O3").Value = 0.1909758294526144
$ws.Range("P3").Value = 0.2199261859093095
$ws.Range("Q3").Value = 0.9762717441363331
$ws.Range("R3").Value = 5.857630464817999
$ws.Range("S3").Value = 0.1391165199792941
$ws.Range("T3").Value = 0.1602054339741137
$ws.Range("G4").Value = 0.1806205
$ws.Range("H4").Value = 0.361241
$ws.Range("I4").Value = 0.7284509268949775
$ws.Range("J4").Value = 0.7284509268949775
$ws.Range("M4").Value = 1.991804666666667
$ws.Range("N4").Value = 5.975414
$ws.Range("O4").Value = 0.07037549633517463
$ws.Range("P4").Value = 0.08104383960437181
$ws.Range("Q4").Value = 0.3597607547956667
$ws.Range("R4").Value = 2.158564528774
$ws.Range("S4").Value = 0.05126509553605205
$ws.Range("T4").Value = 0.05903646007893253
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.1806205
$ws.Range("H5").Value = 0.361241
$ws.Range("I5").Value = 0.7284509268949775
$ws.Range("J5").Value = 0.7284509268949775
$ws.Range("M5").Value = 5.796282
$ws.Range("N5").Value = 11.592564
$ws.Range("O5").Value = 0.2047973024038027
$ws.Range("P5").Value = 0.1572285865748239
$ws.Range("Q5").Value = 1.046927352981
$ws.Range("R5").Value = 4.187709411924
$ws.Range("S5").Value = 0.149184784761641
$ws.Range("T5").Value = 0.1145333096248177
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("G6").Value = 0.1806205
$ws.Range("H6").Value = 0.361241
$ws.Range("I6").Value = 0.7284509268949775
$ws.Range("J6").Value = 0.7284509268949775
$ws.Range("M6").Value = 8.176639
$ws.Range("N6").Value = 24.529917
$ws.Range("O6").Value = 0.2889013353611378
$ws.Range("P6").Value = 0.3326963887115693
$ws.Range("Q6").Value = 1.4768686244995
$ws.Range("R6").Value = 8.861211746996998
$ws.Range("S6").Value = 0.2104504455250176
$ws.Range("T6").Value = 0.2423529927315544
$ws.Range("G7").Value = 0.1806205
$ws.Range("H7").Value = 0.361241
$ws.Range("I7").Value = 0.7284509268949775
$ws.Range("J7").Value = 0.7284509268949775
$ws.Range("M7").Value = 1.552033
$ws.Range("N7").Value = 4.656098999999999
$ws.Range("O7").Value = 0.05483725112782315
$ws.Range("P7").Value = 0.06315012491821921
$ws.Range("Q7").Value = 0.2803289764765
$ws.Range("R7").Value = 1.681973858859
$ws.Range("S7").Value = 0.03994624641243542
$ws.Range("T7").Value = 0.0460017670302104
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.067331
$ws.Range("H8").Value = 0.134662
$ws.Range("I8").Value = 0.2715490731050226
$ws.Range("J8").Value = 0.2715490731050226
$ws.Range("M8").Value = 5.380673
$ws.Range("N8").Value = 10.761346
$ws.Range("O8").Value = 0.1901127853194472
$ws.Range("P8").Value = 0.1459548742817064
$ws.Range("Q8").Value = 0.362286093763
$ws.Range("R8").Value = 1.449144375052
$ws.Range("S8").Value = 0.05162495063891003
$ws.Range("T8").Value = 0.03963391082635746
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.067331
$ws.Range("H9").Value = 0.134662
$ws.Range("I9").Value = 0.2715490731050226
$ws.Range("J9").Value = 0.2715490731050226
$ws.Range("M9").Value = 5.405099333333332
$ws.Range("O9").Value = 0.1909758294526144
$ws.Range("P9").Value = 0.2199261859093095
$ws.Range("Q9").Value = 0.3639307432126666
$ws.Range("R9").Value = 2.183584459276
$ws.Range("S9").Value = 0.05185930947332031
$ws.Range("T9").Value = 0.05972075193519587
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.067331
$ws.Range("H10").Value = 0.134662
$ws.Range("I10").Value = 0.2715490731050226
$ws.Range("J10").Value = 0.2715490731050226
$ws.Range("M10").Value = 1.991804666666667
$ws.Range("N10").Value = 5.975414
$ws.Range("O10").Value = 0.07037549633517463
$ws.Range("P10").Value = 0.08104383960437181
$ws.Range("Q10").Value = 0.1341102000113333
$ws.Range("R10").Value = 0.804661200068
$ws.Range("S10").Value = 0.01911040079912258
$ws.Range("T10").Value = 0.02200737952543928
$ws.Range("D11").Value = "MuSCs"
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.067331
$ws.Range("H11").Value = 0.134662
$ws.Range("I11").Value = 0.2715490731050226
$ws.Range("J11").Value = 0.2715490731050226
$ws.Range("M11").Value = 5.796282
$ws.Range("N11").Value = 11.592564
$ws.Range("O11").Value = 0.2047973024038027
$ws.Range("P11").Value = 0.1572285865748239
$ws.Range("Q11").Value = 0.390269463342
$ws.Range("R11").Value = 1.561077853368
$ws.Range("S11").Value = 0.05561251764216163
$ws.Range("T11").Value = 0.04269527695000623
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.067331
$ws.Range("H12").Value = 0.134662
$ws.Range("I12").Value = 0.2715490731050226
$ws.Range("J12").Value = 0.2715490731050226
$ws.Range("M12").Value = 8.176639
$ws.Range("N12").Value = 24.529917
$ws.Range("O12").Value = 0.2889013353611378
$ws.Range("P12").Value = 0.3326963887115693
$ws.Range("Q12").Value = 0.550541280509
$ws.Range("R12").Value = 3.303247683054
$ws.Range("S12").Value = 0.07845088983612025
$ws.Range("T12").Value = 0.09034339598001494
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.067331
$ws.Range("H13").Value = 0.134662
$ws.Range("I13").Value = 0.2715490731050226
$ws.Range("J13").Value = 0.2715490731050226
$ws.Range("M13").Value = 1.552033
$ws.Range("N13").Value = 4.656098999999999
$ws.Range("O13").Value = 0.05483725112782315
$ws.Range("P13").Value = 0.06315012491821921
$ws.Range("Q13").Value = 0.104499933923
$ws.Range("R13").Value = 0.6269996035379999
$ws.Range("S13").Value = 0.01489100471538773
$ws.Range("T13").Value = 0.01714835788800882

Write-Output "done"
